$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Arttt"
$ws.Range("A14").Value = "Sound"
$ws.Range("A28").Value = "Game"

$ws.Range("A28").Select()
